# Generate Report for Handoff
# Adds a new handed-off file (b013047c-62ce-4ac6-abd6-31ea53f6aaa5.md) as a
# row just above the trailing ".localization-config" row on every sheet.

$wb = $excel.ActiveWorkbook

$newMdBase   = "b013047c-62ce-4ac6-abd6-31ea53f6aaa5.md"
$newMdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/7730bf134b9cbf4346ee10ea1e757ea45e0ebc45/e2e/b013047c-62ce-4ac6-abd6-31ea53f6aaa5.md"
$configDisp  = ".localization-config"
$configUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/7730bf134b9cbf4346ee10ea1e757ea45e0ebc45/.localization-config"
$mdDisp      = "14afe119-6ba8-4e81-a647-21875d2dcef6.md"
$mdUrl       = "https://github.com/OpenLocalizationTest/oltest/blob/7730bf134b9cbf4346ee10ea1e757ea45e0ebc45/e2e/14afe119-6ba8-4e81-a647-21875d2dcef6.md"

# ---------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# push the ".localization-config" summary row from row 3 down to row 4
$ws1.Cells.Item(4,1).Value = $configDisp
$ws1.Cells.Item(4,2).Value = "Not to be localized"
$ws1.Cells.Item(4,3).Value = "Not to be localized"

# write the new handoff row into row 3
$ws1.Cells.Item(3,1).Value = $newMdBase
$ws1.Cells.Item(3,2).Value = "Ready for handoff"
$ws1.Cells.Item(3,3).Value = "Ready for handoff"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdUrl, "", "", $mdDisp)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $newMdUrl, "", "", $newMdBase)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", $configDisp)

# ---------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$zhXlfBase = "b013047c-62ce-4ac6-abd6-31ea53f6aaa5.ebb387106cc1ecce4d5451472042925c952acf69.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d93ad86047a7323c5707e4527287ee0ee6db62a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b013047c-62ce-4ac6-abd6-31ea53f6aaa5.ebb387106cc1ecce4d5451472042925c952acf69.zh-cn.xlf"
$zhOldXlfBase = "14afe119-6ba8-4e81-a647-21875d2dcef6.82e52fd1c53401c582cd9937ded7b19b15574519.zh-cn.xlf"
$zhOldXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8d93ad86047a7323c5707e4527287ee0ee6db62a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/14afe119-6ba8-4e81-a647-21875d2dcef6.82e52fd1c53401c582cd9937ded7b19b15574519.zh-cn.xlf"

# push the ".localization-config" detail row from row 3 down to row 4
$ws2.Cells.Item(4,1).Value = $configDisp
$ws2.Cells.Item(4,2).Value = "Not to be localized"
$ws2.Cells.Item(4,4).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,8).Value = "Ignored"

# write the new handoff detail row into row 3
$ws2.Cells.Item(3,1).Value = $newMdBase
$ws2.Cells.Item(3,2).Value = "Ready for handoff"
$ws2.Cells.Item(3,3).Value = $zhXlfBase
$ws2.Cells.Item(3,4).Value = "2016-03-09 10:34:59"
$ws2.Cells.Item(3,7).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(3,8).Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdUrl, "", "", $mdDisp)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhOldXlfUrl, "", "", $zhOldXlfBase)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $newMdUrl, "", "", $newMdBase)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $zhXlfUrl, "", "", $zhXlfBase)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", $configDisp)

# ---------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$deXlfBase = "b013047c-62ce-4ac6-abd6-31ea53f6aaa5.ebb387106cc1ecce4d5451472042925c952acf69.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a880abce194450bbf1675ec5497450ae22b2c7c6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b013047c-62ce-4ac6-abd6-31ea53f6aaa5.ebb387106cc1ecce4d5451472042925c952acf69.de-de.xlf"
$deOldXlfBase = "14afe119-6ba8-4e81-a647-21875d2dcef6.82e52fd1c53401c582cd9937ded7b19b15574519.de-de.xlf"
$deOldXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a880abce194450bbf1675ec5497450ae22b2c7c6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/14afe119-6ba8-4e81-a647-21875d2dcef6.82e52fd1c53401c582cd9937ded7b19b15574519.de-de.xlf"

# push the ".localization-config" detail row from row 3 down to row 4
$ws3.Cells.Item(4,1).Value = $configDisp
$ws3.Cells.Item(4,2).Value = "Not to be localized"
$ws3.Cells.Item(4,4).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,8).Value = "Ignored"

# write the new handoff detail row into row 3
$ws3.Cells.Item(3,1).Value = $newMdBase
$ws3.Cells.Item(3,2).Value = "Ready for handoff"
$ws3.Cells.Item(3,3).Value = $deXlfBase
$ws3.Cells.Item(3,4).Value = "2016-03-09 10:35:05"
$ws3.Cells.Item(3,7).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(3,8).Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdUrl, "", "", $mdDisp)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deOldXlfUrl, "", "", $deOldXlfBase)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $newMdUrl, "", "", $newMdBase)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $deXlfUrl, "", "", $deXlfBase)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", $configDisp)

Write-Output "Handback report rows generated."
